# Apply the commit "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" to the active workbook.
#
# The only functional change in the target diff is a rewrite of the values in
# column G (header "K", row 1) for data rows 2-71 on Sheet1. Everything else
# in the sheet is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values, keyed by row number (2-71), taken from the
# regenerated save data.
$newK = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 2
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 3
    15 = 1
    16 = 3
    17 = 0
    18 = 3
    19 = 0
    20 = 1
    21 = 1
    22 = 2
    23 = 0
    24 = 2
    25 = 1
    26 = 0
    27 = 3
    28 = 4
    29 = 0
    30 = 1
    31 = 0
    32 = 2
    33 = 2
    34 = 2
    35 = 2
    36 = 0
    37 = 3
    38 = 1
    39 = 1
    40 = 2
    41 = 2
    42 = 3
    43 = 3
    44 = 1
    45 = 2
    46 = 1
    47 = 2
    48 = 1
    49 = 1
    50 = 1
    51 = 4
    52 = 2
    53 = 2
    54 = 0
    55 = 1
    56 = 2
    57 = 0
    58 = 2
    59 = 2
    60 = 2
    61 = 2
    62 = 1
    63 = 0
    64 = 0
    65 = 1
    66 = 1
    67 = 1
    68 = 2
    69 = 2
    70 = 1
    71 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
